$d = $word.ActiveDocument

# --- Step 1: add the new rows to the Attorney table (2nd table in the doc) ---
$t = $d.Tables.Item(2)

$rowsData = @(
    @("[Name]", "Brunson, Barbara"),
    @("[Work]", "(803) 799-0425"),
    @("[Address]", "3614 Landmark Drive"),
    @("", "Suite B"),
    @("", "Columbia, SC 29204"),
    @("", "USA"),
    @("[Law Practice]", "Principal, Law Offices of Barbara E. Brunson"),
    @("[Website]", "http://brunsonlawsc.com/")
)

$blankCells = @()
foreach ($pair in $rowsData) {
    $newRow = $t.Rows.Add()
    if ($pair[0] -ne "") {
        $newRow.Cells.Item(1).Range.Text = $pair[0]
    } else {
        # remember this cell; its alignment mark gets cleared after all rows
        # are created (clearing it earlier would get cloned into later rows
        # via Rows.Add)
        $blankCells += $newRow.Cells.Item(1)
    }
    $newRow.Cells.Item(2).Range.Text = $pair[1]
}

foreach ($cell in $blankCells) {
    $cell.Range.ParagraphFormat.Alignment = 0
}

# --- Step 2: remove the two trailing BodyText paragraphs that used to hold
#             the pipe-delimited text, now superseded by the table rows above ---
$paras = $d.Content.Paragraphs
for ($i = $paras.Count; $i -ge 1; $i--) {
    $p = $paras.Item($i)
    $text = $p.Range.Text
    if ($text.Contains("[Name] | Brunson, Barbara |") -or $text.Contains("|| Suite B |")) {
        $p.Range.Delete()
    }
}
